# Update cryptos list data (prices & 1h volume change %) on sheet1,
# plus swap the FraxShare/Algorand rows (41 <-> 42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text-formatted numbers (e.g. "1.000", "0.07030").
# Force text number-format before assigning so Excel does not coerce these
# into numeric values and strip significant trailing/leading zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.851.36"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.77"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.92"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3635"
$ws.Range("E7").Value = "  -1.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.45"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3287"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.145"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07109"
$ws.Range("E11").Value = "  -2.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.088"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.75"
$ws.Range("E14").Value = "  -3.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.667.07"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.658"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001055"
$ws.Range("E17").Value = "  -3.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06653"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9993"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "79.80"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.967"
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.88"
$ws.Range("E22").Value = "  -5.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.70"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.848.10"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.433"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.438"
$ws.Range("E26").Value = "  -9.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.57"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.74"
$ws.Range("E28").Value = "  -5.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.239"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.850.79"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.18"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.131"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.887"
$ws.Range("E33").Value = "  -9.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08534"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.35"
$ws.Range("E36").Value = "  -7.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.282"
$ws.Range("E37").Value = "  +3.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.243"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02277"
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06099"
$ws.Range("E40").Value = "  -5.53%  "

# Rows 41 and 42 swap content: FraxShare moves to row 42, Algorand moves to row 41.
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2089"
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.349"
$ws.Range("E42").Value = "  -5.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9995"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5984"
$ws.Range("E44").Value = "  -4.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.823"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("E46").Value = "  -4.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5677"
$ws.Range("E47").Value = "  -4.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.19"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.975"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07030"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.206"
$ws.Range("E51").Value = "  +0.78%  "
